# fix: fixed formatting when scrapping floating point numbers
#
# 1) Four "Razon social"/"Nombre Fantasia" entries had a stray comma turned
#    into a period (or dropped) by the scraper - fix the literal text.
# 2) The "Importe" column (H) held amounts as literal text using the
#    Argentine convention (dot thousands separator, comma decimal separator),
#    e.g. "25.000,50". The scraper bug needs them normalised to the plain
#    floating point form "25000.50" (no thousands separator, dot decimal),
#    while staying literal text (not real numbers) in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the four mis-punctuated name/company strings -------------------

$nameFixes = @(
    @{ Cells = @("E31","F31","E41","F41","E42","F42","E56","F56","E81","F81","E90","F90","E134","F134");
       Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA" },
    @{ Cells = @("E72","F72");
       Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH" },
    @{ Cells = @("E80");
       Value = "FERNANDEZ MARIO H. GALLICET OSCAR M" },
    @{ Cells = @("E83");
       Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO" }
)

foreach ($fix in $nameFixes) {
    foreach ($addr in $fix.Cells) {
        $ws.Range($addr).Value = $fix.Value
    }
}

# --- 2) Re-format the "Importe" column (H) text values ---------------------
# Values are stored as literal text (not numbers) formatted like "25.000,50".
# Strip the "." thousands separator and turn the "," decimal separator into
# ".", e.g. "25.000,50" -> "25000.50". Writing the digits straight back with
# Range.Value would make Excel auto-convert the text into a real number
# (and silently drop the trailing zero), so the cell is put into Text mode
# first and reset back to its original (default/"Normal") style afterwards
# so no visible formatting/style change is left behind.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Text
    if ([string]::IsNullOrEmpty($old)) { continue }

    $cell.NumberFormat = "@"
    $cell.Value = $old
    [void]$cell.Replace(".", "")
    [void]$cell.Replace(",", ".")
    $cell.Style = "Normal"
}
